$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2435.3845
$ws.Range("J17").Value = 2435.3845
$ws.Range("L17").Value = 7306.1535
$ws.Range("N17").Value = -7642.1535

$ws.Range("H28").Value = 1308.0952
$ws.Range("I28").Value = 1039
$ws.Range("K28").Value = 1039
$ws.Range("M28").Value = -554

$ws.Range("H40").Value = 5735.727
$ws.Range("I40").Value = 5649.375
$ws.Range("K40").Value = 5649.375
$ws.Range("M40").Value = -5474.375

$ws.Range("H46").Value = 1671198.4
$ws.Range("I46").Value = 5002606
$ws.Range("J46").Value = 5494.5
$ws.Range("K46").Value = 15007818
$ws.Range("L46").Value = 16483.5
$ws.Range("M46").Value = -15007699
$ws.Range("N46").Value = -16721.5

$ws.Range("H60").Value = 1671198.4
$ws.Range("I60").Value = 5002606
$ws.Range("J60").Value = 5494.5
$ws.Range("K60").Value = 15007818
$ws.Range("L60").Value = 16483.5
$ws.Range("M60").Value = -15007334
$ws.Range("N60").Value = -17451.5

$ws.Range("H80").Value = 919.6
$ws.Range("J80").Value = 990.1667
$ws.Range("L80").Value = 2970.5001
$ws.Range("N80").Value = -4966.5001

$ws.Range("H82").Value = 1023.375
$ws.Range("I82").Value = 1128.1428
$ws.Range("J82").Value = 290
$ws.Range("K82").Value = 3384.4284
$ws.Range("L82").Value = 870
$ws.Range("M82").Value = -2978.4284
$ws.Range("N82").Value = -1682

$ws.Range("H83").Value = 919.6
$ws.Range("J83").Value = 990.1667
$ws.Range("L83").Value = 8911.5003
$ws.Range("N83").Value = -18895.5003

$ws.Range("H85").Value = 1023.375
$ws.Range("I85").Value = 1128.1428
$ws.Range("J85").Value = 290
$ws.Range("K85").Value = 3384.4284
$ws.Range("L85").Value = 870
$ws.Range("M85").Value = -1980.4284
$ws.Range("N85").Value = -3678

$ws.Range("H86").Value = 5900
$ws.Range("I86").Value = 5900
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 5900
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -4777
$ws.Range("N86").ClearContents()

$ws.Range("H88").Value = 4946.5
$ws.Range("I88").Value = 1894
$ws.Range("J88").Value = 5964
$ws.Range("K88").Value = 1894
$ws.Range("L88").Value = 5964
$ws.Range("M88").Value = -1488
$ws.Range("N88").Value = -6776

$ws.Range("H89").Value = 5900
$ws.Range("I89").Value = 5900
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 29500
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -23884
$ws.Range("N89").ClearContents()

$ws.Range("H91").Value = 4946.5
$ws.Range("I91").Value = 1894
$ws.Range("J91").Value = 5964
$ws.Range("K91").Value = 1894
$ws.Range("L91").Value = 5964
$ws.Range("M91").Value = -490
$ws.Range("N91").Value = -8772

$ws.Range("H113").Value = 2712
$ws.Range("I113").Value = 2854
$ws.Range("K113").Value = 2854
$ws.Range("M113").Value = 400

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H129").Value = 1788.75
$ws.Range("I129").Value = 1788.75
$ws.Range("K129").Value = 5366.25
$ws.Range("M129").Value = -366.25

$ws.Range("H132").Value = 2358.8462
$ws.Range("I132").Value = 2213.08
$ws.Range("K132").Value = 6639.24
$ws.Range("M132").Value = -4109.24

$ws.Range("H135").Value = 6226.4546
$ws.Range("I135").Value = 1314.8948
$ws.Range("K135").Value = 11834.0532
$ws.Range("M135").Value = -9299.0532

$ws.Range("H138").Value = 2947.311
$ws.Range("I138").Value = 1670.5714
$ws.Range("J138").Value = 3182.5
$ws.Range("K138").Value = 5011.7142
$ws.Range("L138").Value = 9547.5
$ws.Range("M138").Value = 128.2857999999997
$ws.Range("N138").Value = -19827.5

$ws.Range("H141").Value = 3071.7407
$ws.Range("I141").Value = 3071.7407
$ws.Range("K141").Value = 9215.222099999999
$ws.Range("M141").Value = -4035.222099999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3642.8
$ws.Range("I2").Value = 4016.889
$ws.Range("K2").Value = 4016.889
$ws.Range("M2").Value = -3903.889

$ws.Range("H32").Value = 8313.960999999999
$ws.Range("I32").Value = 1105.9736
$ws.Range("K32").Value = 1105.9736
$ws.Range("M32").Value = -818.9736

$ws.Range("H46").Value = 14803.909
$ws.Range("I46").Value = 10233.25
$ws.Range("J46").Value = 17415.715
$ws.Range("K46").Value = 10233.25
$ws.Range("L46").Value = 17415.715
$ws.Range("M46").Value = -9914.25
$ws.Range("N46").Value = -18053.715

$ws.Range("H57").Value = 9000
$ws.Range("I57").Value = 9000
$ws.Range("K57").Value = 9000
$ws.Range("M57").Value = -8516

$ws.Range("H61").Value = 2943.8667
$ws.Range("I61").Value = 2813.6296
$ws.Range("K61").Value = 2813.6296
$ws.Range("M61").Value = -2601.6296

$ws.Range("H74").Value = 2606.8125
$ws.Range("I74").Value = 2673.9333
$ws.Range("K74").Value = 2673.9333
$ws.Range("M74").Value = -1799.9333

$ws.Range("H77").Value = 2606.8125
$ws.Range("I77").Value = 2673.9333
$ws.Range("K77").Value = 13369.6665
$ws.Range("M77").Value = -9001.666500000001

$ws.Range("H115").Value = 50000
$ws.Range("J115").Value = 50000
$ws.Range("L115").Value = 50000
$ws.Range("N115").Value = -53134

$ws.Range("H116").Value = 3642.8
$ws.Range("I116").Value = 4016.889
$ws.Range("K116").Value = 4016.889
$ws.Range("M116").Value = -1722.889

$ws.Range("H119").Value = 93852.71000000001
$ws.Range("J119").Value = 93852.71000000001
$ws.Range("L119").Value = 93852.71000000001
$ws.Range("N119").Value = -103528.71

$ws.Range("H120").Value = 50000
$ws.Range("J120").Value = 50000
$ws.Range("L120").Value = 50000
$ws.Range("N120").Value = -59676

$ws.Range("H131").Value = 89999.5
$ws.Range("J131").Value = 89999.5
$ws.Range("L131").Value = 89999.5
$ws.Range("N131").Value = -100079.5

$ws.Range("H132").Value = 2291.7778
$ws.Range("I132").Value = 2348.75
$ws.Range("J132").Value = 2129
$ws.Range("K132").Value = 7046.25
$ws.Range("L132").Value = 6387
$ws.Range("M132").Value = -4516.25
$ws.Range("N132").Value = -11447

$ws.Range("H136").Value = 2943.8667
$ws.Range("I136").Value = 2813.6296
$ws.Range("K136").Value = 8440.888800000001
$ws.Range("M136").Value = -5890.888800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3642.8
$ws.Range("I3").Value = 4016.889
$ws.Range("K3").Value = 4016.889
$ws.Range("M3").Value = -3902.889

$ws.Range("H107").Value = 3310.3333
$ws.Range("I107").Value = 3199
$ws.Range("K107").Value = 3199
$ws.Range("M107").Value = -1279

$ws.Range("H132").Value = 78168.8
$ws.Range("J132").Value = 78168.8
$ws.Range("L132").Value = 78168.8
$ws.Range("N132").Value = -88288.8

$ws.Range("H134").Value = 3154.88
$ws.Range("I134").Value = 3154.88
$ws.Range("K134").Value = 9464.639999999999
$ws.Range("M134").Value = -6929.639999999999

$ws.Range("H140").Value = 85779.336
$ws.Range("J140").Value = 85779.336
$ws.Range("L140").Value = 85779.336
$ws.Range("N140").Value = -96139.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 84999.5
$ws.Range("J20").Value = 84999.5
$ws.Range("L20").Value = 84999.5
$ws.Range("N20").Value = -85471.5

$ws.Range("H30").Value = 84999.5
$ws.Range("J30").Value = 84999.5
$ws.Range("L30").Value = 84999.5
$ws.Range("N30").Value = -85181.5

$ws.Range("H58").Value = 10401.833
$ws.Range("I58").Value = 11482.2
$ws.Range("K58").Value = 11482.2
$ws.Range("M58").Value = -11279.2

$ws.Range("H110").Value = 25000
$ws.Range("J110").Value = 25000
$ws.Range("L110").Value = 25000
$ws.Range("N110").Value = -33180

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H112").Value = 25000
$ws.Range("J112").Value = 25000
$ws.Range("L112").Value = 25000
$ws.Range("N112").Value = -27954

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

$ws.Range("H119").Value = 80000
$ws.Range("J119").Value = 80000
$ws.Range("L119").Value = 80000
$ws.Range("N119").Value = -89676

$ws.Range("H120").Value = 88208.336
$ws.Range("J120").Value = 88208.336
$ws.Range("L120").Value = 88208.336
$ws.Range("N120").Value = -95466.336

$ws.Range("H121").Value = 38085.8
$ws.Range("I121").Value = 22296
$ws.Range("J121").Value = 42033.25
$ws.Range("K121").Value = 22296
$ws.Range("L121").Value = 42033.25
$ws.Range("M121").Value = -20986
$ws.Range("N121").Value = -44653.25

$ws.Range("H123").Value = 89998.5
$ws.Range("J123").Value = 89998.5
$ws.Range("L123").Value = 89998.5
$ws.Range("N123").Value = -99798.5

$ws.Range("H124").Value = 74646
$ws.Range("J124").Value = 68996
$ws.Range("L124").Value = 68996
$ws.Range("N124").Value = -73906

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H127").Value = 80709
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").ClearContents()

$ws.Range("H128").Value = 84999.5
$ws.Range("J128").Value = 84999.5
$ws.Range("L128").Value = 84999.5
$ws.Range("N128").Value = -94959.5

$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()

$ws.Range("H130").Value = 83677
$ws.Range("J130").Value = 96666.336
$ws.Range("L130").Value = 96666.336
$ws.Range("N130").Value = -106706.336

$ws.Range("H131").Value = 64879
$ws.Range("I131").Value = 34644.5
$ws.Range("K131").Value = 34644.5
$ws.Range("M131").Value = -29604.5

$ws.Range("H132").Value = 7305.1333
$ws.Range("I132").Value = 6894.5386
$ws.Range("K132").Value = 20683.6158
$ws.Range("M132").Value = -18153.6158

$ws.Range("H133").Value = 60755.715
$ws.Range("J133").Value = 59999
$ws.Range("L133").Value = 59999
$ws.Range("N133").Value = -65059

$ws.Range("H134").Value = 2573.3333
$ws.Range("I134").Value = 1480
$ws.Range("K134").Value = 4440
$ws.Range("M134").Value = -1905

$ws.Range("H135").Value = 85472.28999999999
$ws.Range("J135").Value = 87519.60000000001
$ws.Range("L135").Value = 87519.60000000001
$ws.Range("N135").Value = -97659.60000000001

$ws.Range("H136").Value = 10401.833
$ws.Range("I136").Value = 11482.2
$ws.Range("K136").Value = 34446.60000000001
$ws.Range("M136").Value = -31896.60000000001

$ws.Range("H137").Value = 86365
$ws.Range("J137").Value = 89962.5
$ws.Range("L137").Value = 89962.5
$ws.Range("N137").Value = -100162.5

$ws.Range("H138").Value = 89996.5
$ws.Range("J138").Value = 89996.5
$ws.Range("L138").Value = 89996.5
$ws.Range("N138").Value = -100276.5

$ws.Range("H139").Value = 42531.91
$ws.Range("J139").Value = 89998.25
$ws.Range("L139").Value = 89998.25
$ws.Range("N139").Value = -100278.25

$ws.Range("H140").Value = 87998
$ws.Range("J140").Value = 87998
$ws.Range("L140").Value = 87998
$ws.Range("N140").Value = -98358

$ws.Range("H141").Value = 73374.375
$ws.Range("J141").Value = 87499.336
$ws.Range("L141").Value = 87499.336
$ws.Range("N141").Value = -97859.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value = 27539.447
$ws.Range("I46").Value = 1476.1904
$ws.Range("J46").Value = 59735.234
$ws.Range("K46").Value = 4428.5712
$ws.Range("L46").Value = 179205.702
$ws.Range("M46").Value = -4337.5712
$ws.Range("N46").Value = -179387.702

$ws.Range("H80").Value = 5000
$ws.Range("I80").Value = 5000
$ws.Range("K80").Value = 15000
$ws.Range("M80").Value = -14064

$ws.Range("H83").Value = 5000
$ws.Range("I83").Value = 5000
$ws.Range("K83").Value = 45000
$ws.Range("M83").Value = -40320

$ws.Range("H121").Value = 504.57144
$ws.Range("I121").Value = 259.8
$ws.Range("K121").Value = 779.4000000000001
$ws.Range("M121").Value = 530.5999999999999

$ws.Range("H122").Value = 2677.8823
$ws.Range("J122").Value = 4065.1
$ws.Range("L122").Value = 36585.9
$ws.Range("N122").Value = -41485.9

$ws.Range("H126").Value = 1997.5
$ws.Range("I126").Value = 1997.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 5992.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1052.5
$ws.Range("N126").ClearContents()

$ws.Range("H132").Value = 1377.625
$ws.Range("I132").Value = 1377.625
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 12398.625
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -9868.625
$ws.Range("N132").ClearContents()

$ws.Range("H134").Value = 6944.25
$ws.Range("I134").Value = 2592.3333
$ws.Range("J134").Value = 20000
$ws.Range("K134").Value = 7776.999899999999
$ws.Range("L134").Value = 60000
$ws.Range("M134").Value = -2706.999899999999
$ws.Range("N134").Value = -70140

$ws.Range("H140").Value = 4297.25
$ws.Range("I140").Value = 3443.4614
$ws.Range("K140").Value = 10330.3842
$ws.Range("M140").Value = -5150.3842

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5913.4
$ws.Range("J70").Value = 5889.5
$ws.Range("L70").Value = 5889.5
$ws.Range("N70").Value = -6429.5

$ws.Range("H73").Value = 5913.4
$ws.Range("J73").Value = 5889.5
$ws.Range("L73").Value = 5889.5
$ws.Range("N73").Value = -7761.5

$ws.Range("H80").Value = 4843.3335
$ws.Range("I80").Value = 5820
$ws.Range("J80").Value = 2890
$ws.Range("K80").Value = 5820
$ws.Range("L80").Value = 2890
$ws.Range("M80").Value = -4822
$ws.Range("N80").Value = -4886

$ws.Range("H83").Value = 4843.3335
$ws.Range("I83").Value = 5820
$ws.Range("J83").Value = 2890
$ws.Range("K83").Value = 29100
$ws.Range("L83").Value = 14450
$ws.Range("M83").Value = -24108
$ws.Range("N83").Value = -24434

$ws.Range("H102").Value = 2577.5
$ws.Range("I102").Value = 2635.0667
$ws.Range("J102").Value = 1714
$ws.Range("K102").Value = 2635.0667
$ws.Range("L102").Value = 1714
$ws.Range("M102").Value = -1013.0667
$ws.Range("N102").Value = -4958

$ws.Range("H111").Value = 37813.5
$ws.Range("J111").Value = 37813.5
$ws.Range("L111").Value = 37813.5
$ws.Range("N111").Value = -43947.5

$ws.Range("H114").Value = 69000
$ws.Range("J114").Value = 69000
$ws.Range("L114").Value = 69000
$ws.Range("N114").Value = -77678

$ws.Range("H118").Value = 20000
$ws.Range("J118").Value = 20000
$ws.Range("L118").Value = 20000
$ws.Range("N118").Value = -23314

$ws.Range("H119").Value = 80000
$ws.Range("J119").Value = 80000
$ws.Range("L119").Value = 80000
$ws.Range("N119").Value = -89676

$ws.Range("H120").Value = 37000
$ws.Range("J120").Value = 37000
$ws.Range("L120").Value = 37000
$ws.Range("N120").Value = -46676

$ws.Range("H121").Value = 45000
$ws.Range("J121").Value = 45000
$ws.Range("L121").Value = 45000
$ws.Range("N121").Value = -48494

$ws.Range("H122").Value = 2656.1667
$ws.Range("I122").Value = 1984.25
$ws.Range("J122").Value = 4000
$ws.Range("K122").Value = 5952.75
$ws.Range("L122").Value = 12000
$ws.Range("M122").Value = -3502.75
$ws.Range("N122").Value = -16900

$ws.Range("H124").Value = 93642.14
$ws.Range("J124").Value = 93642.14
$ws.Range("L124").Value = 93642.14
$ws.Range("N124").Value = -103462.14

$ws.Range("H128").Value = 18999
$ws.Range("J128").Value = 18999
$ws.Range("L128").Value = 18999
$ws.Range("N128").Value = -28959

$ws.Range("H130").Value = 64990
$ws.Range("J130").Value = 64990
$ws.Range("L130").Value = 64990
$ws.Range("N130").Value = -75030

$ws.Range("H132").Value = 4397.8726
$ws.Range("I132").Value = 4302.41
$ws.Range("J132").Value = 4863.25
$ws.Range("K132").Value = 12907.23
$ws.Range("L132").Value = 14589.75
$ws.Range("M132").Value = -10377.23
$ws.Range("N132").Value = -19649.75

$ws.Range("H133").Value = 82259.5
$ws.Range("J133").Value = 82259.5
$ws.Range("L133").Value = 82259.5
$ws.Range("N133").Value = -92379.5

$ws.Range("H135").Value = 87009.60000000001
$ws.Range("J135").Value = 87009.60000000001
$ws.Range("L135").Value = 87009.60000000001
$ws.Range("N135").Value = -97149.60000000001

$ws.Range("H137").Value = 89780
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 89780
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 89780
$ws.Range("M137").ClearContents()
$ws.Range("N137").Value = -99980

$ws.Range("H138").Value = 104500
$ws.Range("J138").Value = 104500
$ws.Range("L138").Value = 104500
$ws.Range("N138").Value = -114780

$ws.Range("H139").Value = 82962.25
$ws.Range("J139").Value = 82962.25
$ws.Range("L139").Value = 82962.25
$ws.Range("N139").Value = -93242.25

$ws.Range("H140").Value = 79780
$ws.Range("J140").Value = 79780
$ws.Range("L140").Value = 79780
$ws.Range("N140").Value = -90140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H7").Value = 24581.084
$ws.Range("I7").Value = 24997.637
$ws.Range("K7").Value = 24997.637
$ws.Range("M7").Value = -24885.637

$ws.Range("H22").Value = 1977.9
$ws.Range("J22").Value = 1600.6
$ws.Range("L22").Value = 1600.6
$ws.Range("N22").Value = -2190.6

$ws.Range("H27").Value = 1977.9
$ws.Range("J27").Value = 1600.6
$ws.Range("L27").Value = 1600.6
$ws.Range("N27").Value = -1814.6

$ws.Range("H40").Value = 2117.2
$ws.Range("I40").Value = 2117.2
$ws.Range("K40").Value = 2117.2
$ws.Range("M40").Value = -1981.2

$ws.Range("H61").Value = 3428
$ws.Range("I61").Value = 3347.5
$ws.Range("K61").Value = 3347.5
$ws.Range("M61").Value = -3145.5

$ws.Range("H100").Value = 3245.6365
$ws.Range("I100").Value = 1938.9166
$ws.Range("J100").Value = 4813.7
$ws.Range("K100").Value = 1938.9166
$ws.Range("L100").Value = 4813.7
$ws.Range("M100").Value = -1397.9166
$ws.Range("N100").Value = -5895.7

$ws.Range("H113").Value = 3428
$ws.Range("I113").Value = 3347.5
$ws.Range("K113").Value = 3347.5
$ws.Range("M113").Value = -1177.5

$ws.Range("H122").Value = 12205.667
$ws.Range("J122").Value = 2979
$ws.Range("L122").Value = 8937
$ws.Range("N122").Value = -13837

$ws.Range("H126").Value = 24581.084
$ws.Range("I126").Value = 24997.637
$ws.Range("K126").Value = 74992.91099999999
$ws.Range("M126").Value = -72522.91099999999

$ws.Range("H132").Value = 2050.6099
$ws.Range("I132").Value = 1766.0625
$ws.Range("K132").Value = 5298.1875
$ws.Range("M132").Value = -2768.1875

$ws.Range("H136").Value = 1508
$ws.Range("I136").Value = 1417.909
$ws.Range("K136").Value = 4253.727000000001
$ws.Range("M136").Value = -1703.727000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 27654.889
$ws.Range("I45").Value = 33079.332
$ws.Range("K45").Value = 33079.332
$ws.Range("M45").Value = -32588.332

$ws.Range("H113").Value = 1062.7142
$ws.Range("I113").Value = 1071.36
$ws.Range("K113").Value = 3214.08
$ws.Range("M113").Value = -1044.08

$ws.Range("H122").Value = 2134.516
$ws.Range("I122").Value = 2003.5
$ws.Range("J122").Value = 2583.7144
$ws.Range("K122").Value = 6010.5
$ws.Range("L122").Value = 7751.1432
$ws.Range("M122").Value = -3560.5
$ws.Range("N122").Value = -12651.1432

$ws.Range("H126").Value = 2051.8333
$ws.Range("I126").Value = 2051.8333
$ws.Range("K126").Value = 6155.499899999999
$ws.Range("M126").Value = -3685.499899999999

$ws.Range("H132").Value = 4071.697
$ws.Range("I132").Value = 3461.926
$ws.Range("J132").Value = 6815.6665
$ws.Range("K132").Value = 10385.778
$ws.Range("L132").Value = 20446.9995
$ws.Range("M132").Value = -7855.778
$ws.Range("N132").Value = -25506.9995

$ws.Range("H136").Value = 3641.3333
$ws.Range("I136").Value = 3641.3333
$ws.Range("K136").Value = 10923.9999
$ws.Range("M136").Value = -8373.999899999999
